$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.788799643516541
$ws.Range("B1").Value = 1.982946276664734
$ws.Range("C1").Value = 2.372143268585205
$ws.Range("D1").Value = 3.682628154754639
$ws.Range("E1").Value = 1.368188261985779
